# Weekly update ("Fruta / hortaliza, semanal"): a new week's worth of
# Cereza records is inserted at the top of the existing data block
# (rows 289-293), pushing the previously-last week's rows (old 289-301)
# down to 294-306.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data rows 289:301 down by 5 rows, creating 5 blank
# rows at 289:293 for the new data.
$ws.Rows("289:293").Insert()

# Common (unchanging) header values for every Cereza row in this block.
$mercadoId = 9
$mercado   = "Vega Central Mapocho de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$tipo      = "Fruta"
$productoId = 100103
$producto  = "Frutos de hueso (carozo)"
$categoriaId = 100103001
$categoria = "Cereza"

$newRows = @(
    @{ Row=289; Fecha=44578; Variedad="Santina";     Calidad="Especial"; Volumen=300; PMin=6000; PMax=6000; PProm=6000; Unidad="$/bandeja 10 kilos"; Origen="Región de O'Higgins"; PKg=600; KgUnidad=10 },
    @{ Row=290; Fecha=44578; Variedad="Santina";     Calidad="Primera";  Volumen=330; PMin=5000; PMax=5000; PProm=5000; Unidad="$/bandeja 10 kilos"; Origen="Región de O'Higgins"; PKg=500; KgUnidad=10 },
    @{ Row=291; Fecha=44578; Variedad="Santina";     Calidad="Segunda";  Volumen=380; PMin=4000; PMax=4000; PProm=4000; Unidad="$/bandeja 10 kilos"; Origen="Región de O'Higgins"; PKg=400; KgUnidad=10 },
    @{ Row=292; Fecha=44578; Variedad="Sweet Heart"; Calidad="Especial"; Volumen=330; PMin=5500; PMax=5500; PProm=5500; Unidad="$/bandeja 10 kilos"; Origen="Región Metropolitana"; PKg=550; KgUnidad=10 },
    @{ Row=293; Fecha=44578; Variedad="Sweet Heart"; Calidad="Primera";  Volumen=280; PMin=4500; PMax=4500; PProm=4500; Unidad="$/bandeja 10 kilos"; Origen="Región Metropolitana"; PKg=450; KgUnidad=10 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $r.Fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $productoId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $r.Variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $r.Unidad
    $ws.Cells.Item($row, 18).Value = $r.Origen
    $ws.Cells.Item($row, 19).Value = $r.PKg
    $ws.Cells.Item($row, 20).Value = $r.KgUnidad
}
